# ngohrms MenuList.xlsx - "hr portion added"
# Adds a new "Fixed Asset" worksheet (after "HR") populated with the Fixed
# Asset / Vehicle Management menu rows, and makes it the active sheet.

$wb = $excel.ActiveWorkbook

# --- add the new worksheet after the last existing sheet (HR) ----------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Fixed Asset"

# --- header row ----------------------------------------------------------
$ws.Range("A1").Value = "Menu Name"
$ws.Range("B1").Value = "Is Sub Menu"
$ws.Range("C1").Value = "Sub Menu Name"
$ws.Range("D1").Value = "URL"
$ws.Range("E1").Value = "Icon"
$ws.Range("F1").Value = "Ordering"

# --- menu rows -------------------------------------------------------------
$ws.Range("A2").Value = "Fixed Asset"
$ws.Range("B2").Value = "No"
$ws.Range("D2").Value = "/fixed-asset/"
$ws.Range("E2").Value = "icons-Box-Open"

$ws.Range("A3").Value = "Fixed Asset Import"
$ws.Range("B3").Value = "No"
$ws.Range("D3").Value = "/fixed-asset/import/"
$ws.Range("E3").Value = "icons-Box-Open"

$ws.Range("A4").Value = "Maintenance Request"
$ws.Range("B4").Value = "No"
$ws.Range("D4").Value = "/fixed-asset/request-list/"
$ws.Range("E4").Value = "icons-Box-Open"

$ws.Range("A5").Value = "Maintenance Pending "
$ws.Range("B5").Value = "No"
$ws.Range("D5").Value = "/fixed-asset/maintenance-pending-list/"
$ws.Range("E5").Value = "icons-Box-Open"

$ws.Range("A6").Value = "Vehicle Menagement"
$ws.Range("B6").Value = "No"
$ws.Range("D6").Value = "/vehicle-mangement/vehicle/"
$ws.Range("E6").Value = "icons-Box-Open"

# --- Ordering column: "1".."4" stored as TEXT (matches the other sheets,
# where these numeric-looking Ordering values are plain shared-string text,
# not numbers). Route them through a throw-away formula + paste-values so
# no incidental number-format/quote-prefix style gets attached to the cell.
$ws.Range("Z1").Formula = '="1"'
$ws.Range("Z2").Formula = '="2"'
$ws.Range("Z3").Formula = '="3"'
$ws.Range("Z4").Formula = '="4"'
$ws.Range("Z1:Z4").Copy()
$ws.Range("F2:F5").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("Z1:Z4").ClearContents()

# Last row's Ordering is a genuine number, left-aligned.
$ws.Range("F6").Value = 5
$ws.Range("F6").HorizontalAlignment = -4131  # xlLeft

# --- view state: selection + make this the active/selected tab -----------
$ws.Range("E12").Select()
$ws.Activate()
